$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for columns B,D,E,F,G,H,I,L,M,N,O across rows 2-25
$cols = @("B", "D", "E", "F", "G", "H", "I", "L", "M", "N", "O")

$data = @(
    @(13.95311048416221, 3.021000841075117, 21.84986980018313, 18.09084327943867, 19.23381989760912, 11.16362239624939, 23.45497124202965, 8.613128260836604, 13.99936275408525, 18.60532215384137, 15.77549662627584),
    @(13.71403913183276, 2.975564266117376, 21.79147934319596, 17.89144133139383, 18.76531818884532, 11.16504833229945, 23.58991053081361, 8.536192770116221, 13.90019448626649, 18.60599224046507, 15.69982921135823),
    @(13.5666107040239, 2.956068068504205, 21.75708217808031, 17.77369449657812, 18.47905491308302, 11.16822077166237, 23.6771061985449, 8.4896439178164, 13.84043773660549, 18.60847909058312, 15.65765240014532),
    @(13.50644395588689, 2.948216745534363, 21.74343976721936, 17.72694961915625, 18.36296147266134, 11.17009170751185, 23.71373402415408, 8.470864948085644, 13.81639011250459, 18.6100167732415, 15.64155704409685),
    @(13.49645009540426, 2.946918931229818, 21.74119734238752, 17.71926411505339, 18.34372435535732, 11.17043730433871, 23.71988226190565, 8.467758672476803, 13.81241590916713, 18.61030383667921, 15.63895078301438),
    @(13.56579953651994, 2.955961793097467, 21.75689666296382, 17.77305899017524, 18.47748668377851, 11.16824366235646, 23.67759573713476, 8.489389866819122, 13.8401121675842, 18.60849770247574, 15.65743089242551),
    @(13.87085371586262, 3.001220925319708, 21.82944021727103, 18.02115321061076, 19.07211418886468, 11.16363775394355, 23.50059882452279, 8.58646777515612, 13.96494582556428, 18.60512391134727, 15.74852699667855),
    @(14.46092606348855, 3.139583061728506, 21.98284413006223, 18.54211733232056, 20.24038673859481, 11.17279510262513, 23.18782796008752, 8.781569262148533, 14.21788786823367, 18.61486762332304, 15.9604157194541),
    @(14.88551064521858, 3.235169265904108, 22.10180367399632, 18.94207153420929, 21.08874873832353, 11.1905485494308, 22.97876666886652, 8.926789561515074, 14.40753463456773, 18.63185090261868, 16.13526642250771),
    @(15.07592991137344, 3.277237963557564, 22.15716022924797, 19.12692426962258, 21.47033571344265, 11.20099963342469, 22.88812045659785, 8.993046700161129, 14.4943867996847, 18.64167637472136, 16.21871156477589),
    @(15.14758236025306, 3.292957443738364, 22.17829039001439, 19.19727180057389, 21.61404256316742, 11.2052965228015, 22.85443297279715, 9.018146650072774, 14.52733770228645, 18.64569584826069, 16.2508475574561),
    @(15.13217193924537, 3.289581474480722, 22.17373233864502, 19.18210685597826, 21.58313048622269, 11.20435605940748, 22.86165982930303, 9.012740780564172, 14.5202387288105, 18.64481694588751, 16.24390300672702),
    @(15.0818342744362, 3.278535486014127, 22.15889531821779, 19.13270523714566, 21.48217535182678, 11.20134635503406, 22.88533619003279, 8.995111595957198, 14.49709657429163, 18.6420010897191, 16.22134479072107),
    @(15.0509399025362, 3.271741798610493, 22.14982875292285, 19.10248851784275, 21.42022951335036, 11.19954693864119, 22.89992168904483, 8.984313947576084, 14.48292875130738, 18.64031511087488, 16.20759644066221),
    @(14.87300596565914, 3.232390854718187, 22.09821006284335, 18.93004349617954, 21.06370920002803, 11.18991315342982, 22.98478007031924, 8.922461869090286, 14.4018684791182, 18.63125076941362, 16.12988960108779),
    @(14.7631050280164, 3.207882527515803, 22.06685405301539, 18.82494927824717, 20.84376113043674, 11.18460993467997, 23.03797756606218, 8.884554656541461, 14.35227400442456, 18.6262256673511, 16.08320207629615),
    @(14.69963987370576, 3.193653323786003, 22.04893638705346, 18.76478015178456, 20.71685387177315, 11.18178335774902, 23.06899500983432, 8.862770940028094, 14.32380455048188, 18.62353319107502, 16.05671770637045),
    @(14.67811017905817, 3.188813018627096, 22.04289027292534, 18.74445787964705, 20.67382206100027, 11.18086480985804, 23.07956913903784, 8.855399259669355, 14.31417555628763, 18.6226556363228, 16.04781466421431),
    @(14.77483082687932, 3.210505268320492, 22.07017987161132, 18.83610844978767, 20.86721752475292, 11.18515133499167, 23.03227118853972, 8.888588058141691, 14.35754778447407, 18.62674014587884, 16.08813401971349),
    @(15.09663250784249, 3.281785746321545, 22.16324884147558, 19.14720681336517, 21.5118510917812, 11.20222118915119, 22.87836456720712, 9.00028959385179, 14.50389247942187, 18.64282009185377, 16.22795630751778),
    @(15.30426830429993, 3.327138470334056, 22.22504801514036, 19.3525224055754, 21.92847158896611, 11.21535372501888, 22.78149682400817, 9.073341811961114, 14.59988781362466, 18.65506923816456, 16.3224570261643),
    @(15.19371443058325, 3.303048107506203, 22.19197914248623, 19.24278236853256, 21.70659441053733, 11.20816461629941, 22.83285751719215, 9.034354011922805, 14.54862836562617, 18.64837350405941, 16.27174316377314),
    @(14.76953046561357, 3.209319959693147, 22.06867592772347, 18.83106260125896, 20.85661429682824, 11.18490587524013, 23.0348496913238, 8.886764525828719, 14.35516337269316, 18.62650693775528, 16.08590317449159),
    @(14.30258727726826, 3.103182902902393, 21.94020608812513, 18.39788136063837, 19.92530757286121, 11.16837583733166, 23.26878644044808, 8.728385899234098, 14.14870658843297, 18.61049514527085, 15.89964164860299)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $startRow + $r
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $cellRef = "$($cols[$c])$rowNum"
        $ws.Range($cellRef).Value = $rowValues[$c]
    }
}
